# Natmi following Dr Hou advice
# Update LR-pair stats (Rspo2-Lgr4) for the recomputed receptor-expressing
# cell counts, and add a new "sCs" target-cluster row (row 5) while the
# former row 4 (which used to target "sCs") now targets the new "M2"
# cluster.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (target cluster: ECs) ---
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.876175666666667
$ws.Range("H2").Value = 5.628527
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.574348666666667
$ws.Range("N2").Value = 4.723046
$ws.Range("O2").Value = 0.0919647296207302
$ws.Range("P2").Value = 0.107609035279699
$ws.Range("Q2").Value = 2.953754659249111
$ws.Range("R2").Value = 26.583791933242
$ws.Range("S2").Value = 0.0919647296207302
$ws.Range("T2").Value = 0.107609035279699

# --- Row 3 (target cluster: FAPs) ---
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.876175666666667
$ws.Range("H3").Value = 5.628527
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 7.791016
$ws.Range("N3").Value = 23.373048
$ws.Range("O3").Value = 0.4551080043963893
$ws.Range("P3").Value = 0.532527345028208
$ws.Range("Q3").Value = 14.61731463781067
$ws.Range("R3").Value = 131.555831740296
$ws.Range("S3").Value = 0.4551080043963893
$ws.Range("T3").Value = 0.532527345028208

# --- Row 4 (target cluster was sCs, now M2) ---
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.876175666666667
$ws.Range("H4").Value = 5.628527
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.2873323333333334
$ws.Range("N4").Value = 0.861997
$ws.Range("O4").Value = 0.01678436353126363
$ws.Range("P4").Value = 0.01963958546751285
$ws.Range("Q4").Value = 0.5390859320465556
$ws.Range("R4").Value = 4.851773388419
$ws.Range("S4").Value = 0.01678436353126363
$ws.Range("T4").Value = 0.01963958546751285

# --- Row 5 (new row, target cluster: sCs) ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rspo2"
$ws.Range("C5").Value = "Lgr4"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.876175666666667
$ws.Range("H5").Value = 5.628527
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.4663515
$ws.Range("N5").Value = 14.932703
$ws.Range("O5").Value = 0.436142902451617
$ws.Range("P5").Value = 0.3402240342245802
$ws.Range("Q5").Value = 14.00818700308017
$ws.Range("R5").Value = 84.049122018481
$ws.Range("S5").Value = 0.436142902451617
$ws.Range("T5").Value = 0.3402240342245802
